# Apply row-shuffle + field edits to the Artfynd sheet (rows 2-9)
# The underlying rows keep their shared attributes (species, location, dates, etc.)
# but the row-specific fields (Id, TaxonId, Ost/Nord coords, Aktivitet/Metod, comments, ...)
# get reassigned following the recorded permutation, and the Starttid/Sluttid columns are
# dropped for every record row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 112194228
$ws.Range("B2").Value = 89423
$ws.Range("E2").Value = 5432
$ws.Range("F2").Value = "Granticka"
$ws.Range("G2").Value = "Porodaedalea chrysoloma"
$ws.Range("H2").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("I2").Value = ""
$ws.Range("I2").Font.Bold = $ws.Range("I2").Font.Bold
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = ""
$ws.Range("Q2").Value = 511193
$ws.Range("R2").Value = 7079842
$ws.Range("Z2").Value = ""
$ws.Range("AB2").Value = ""

# Row 3
$ws.Range("A3").Value = 112194226
$ws.Range("K3").Value = ""
$ws.Range("L3").Value = ""
$ws.Range("M3").Value = ""
$ws.Range("N3").Value = ""
$ws.Range("Q3").Value = 511177
$ws.Range("R3").Value = 7079788
$ws.Range("Z3").Value = ""
$ws.Range("AB3").Value = ""
$ws.Range("AC3").Value = ""

# Row 4
$ws.Range("A4").Value = 112194221
$ws.Range("Q4").Value = 511121
$ws.Range("R4").Value = 7080073
$ws.Range("Z4").Value = ""
$ws.Range("AB4").Value = ""

# Row 5
$ws.Range("A5").Value = 112194223
$ws.Range("Q5").Value = 510986
$ws.Range("R5").Value = 7079917
$ws.Range("Z5").Value = ""
$ws.Range("AB5").Value = ""
$ws.Range("AC5").Value = "ringhack"

# Row 6
$ws.Range("A6").Value = 112194222
$ws.Range("I6").Value = "1"
$ws.Range("K6").Value = ""
$ws.Range("K6").Font.Bold = $ws.Range("K6").Font.Bold
$ws.Range("L6").Value = ""
$ws.Range("L6").Font.Bold = $ws.Range("L6").Font.Bold
$ws.Range("M6").Value = "födosökande"
$ws.Range("N6").Value = "observerad"
$ws.Range("Q6").Value = 511065
$ws.Range("R6").Value = 7080083
$ws.Range("Z6").Value = ""
$ws.Range("AB6").Value = ""

# Row 7
$ws.Range("A7").Value = 112194219
$ws.Range("B7").Value = 56398
$ws.Range("E7").Value = 100109
$ws.Range("F7").Value = "Tretåig hackspett"
$ws.Range("G7").Value = "Picoides tridactylus"
$ws.Range("H7").Value = "(Linnaeus, 1758)"
$ws.Range("K7").Value = ""
$ws.Range("K7").Font.Bold = $ws.Range("K7").Font.Bold
$ws.Range("L7").Value = ""
$ws.Range("L7").Font.Bold = $ws.Range("L7").Font.Bold
$ws.Range("M7").Value = ""
$ws.Range("M7").Font.Bold = $ws.Range("M7").Font.Bold
$ws.Range("N7").Value = ""
$ws.Range("N7").Font.Bold = $ws.Range("N7").Font.Bold
$ws.Range("Q7").Value = 511201
$ws.Range("R7").Value = 7080015
$ws.Range("Z7").Value = ""
$ws.Range("AB7").Value = ""
$ws.Range("AC7").Value = "ringhack äldre"

# Row 8
$ws.Range("A8").Value = 112194225
$ws.Range("K8").Value = ""
$ws.Range("L8").Value = ""
$ws.Range("M8").Value = ""
$ws.Range("N8").Value = ""
$ws.Range("Q8").Value = 511026
$ws.Range("R8").Value = 7079794
$ws.Range("Z8").Value = ""
$ws.Range("AB8").Value = ""
$ws.Range("AC8").Value = ""

# Row 9
$ws.Range("A9").Value = 112194220
$ws.Range("K9").Value = ""
$ws.Range("K9").Font.Bold = $ws.Range("K9").Font.Bold
$ws.Range("L9").Value = ""
$ws.Range("L9").Font.Bold = $ws.Range("L9").Font.Bold
$ws.Range("M9").Value = ""
$ws.Range("M9").Font.Bold = $ws.Range("M9").Font.Bold
$ws.Range("N9").Value = ""
$ws.Range("N9").Font.Bold = $ws.Range("N9").Font.Bold
$ws.Range("Q9").Value = 511163
$ws.Range("R9").Value = 7080066
$ws.Range("Z9").Value = ""
$ws.Range("AB9").Value = ""
$ws.Range("AC9").Value = "ringhack äldre"
